# 035 Week 19/20 update
# Fills in the WK20 scores (Sheet1, column W) and the WK19 scores
# (THURSDAY SINGLES, column T) that were posted for weeks 19/20, and bumps
# the HANDICAPS sheet's week-9 game count (ALBIE GILLESPIE) by one.
# All SUM/COUNTIF formulas downstream recalculate automatically.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: WK20 scores (column W) -------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("W12").Value = 33
$ws1.Range("W13").Value = 28
$ws1.Range("W14").Value = 31
$ws1.Range("W16").Value = 33
$ws1.Range("W17").Value = 34
$ws1.Range("W18").Value = 33
$ws1.Range("W20").Value = 31
$ws1.Range("W21").Value = 36
$ws1.Range("W22").Value = 32
$ws1.Range("W23").Value = 34
$ws1.Range("W24").Value = 31
$ws1.Range("W25").Value = 30
$ws1.Range("W26").Value = 30

# W29 is highlighted red in the source sheet (matches the style already
# used for other standout scores in that row), same as the original file.
$ws1.Range("W29").Value = 36
$ws1.Range("W29").Font.Color = 255

# ---- THURSDAY SINGLES: WK19 scores (column T) ----------------------------
$ws2 = $wb.Worksheets.Item("THURSDAY SINGLES")

# T6 keeps the red-highlighted font used for this cell in the source file.
$ws2.Range("T6").Value = 34
$ws2.Range("T6").Font.Color = 255

$ws2.Range("T9").Value = 22
$ws2.Range("T15").Value = 25

# ---- HANDICAPS: ALBIE GILLESPIE's week-9 game count ----------------------
$ws3 = $wb.Worksheets.Item("HANDICAPS")
$ws3.Range("B10").Value = 7
